# Trade #110 closed at 2026-02-18 00:38:05 - unknown UNKNOWN +0.000%
#
# This script applies the bookkeeping updates that result from:
#   1) Closing the open HighProbConvergence trade (row 139 on "All Trades",
#      row 14 on "HighProbConvergence") with an early exit.
#   2) Recomputing the Summary and Strategy Status roll-up numbers.
#   3) Appending two brand-new OPEN trades (one HighProbConvergence, one
#      MarketMaking) to the "All Trades" sheet and to each strategy's own
#      sheet.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1) Summary sheet
# -----------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.1     # Current Capital
$summary.Range("B4").Value = 0.21       # Total P&L $
$summary.Range("B6").Value = 138        # Total Trades
$summary.Range("B8").Value = 49         # Losing Trades
$summary.Range("B9").Value = 46.38      # Win Rate %

# -----------------------------------------------------------------------
# 2) Strategy Status sheet - HighProbConvergence row (row 3)
# -----------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C3").Value = 100.4       # Capital
$status.Range("D3").Value = 13          # Trades
$status.Range("E3").Value = 0.41        # P&L $
$status.Range("F3").Value = 0.4         # P&L %
$status.Range("G3").Value = 76.92       # Win Rate %

# -----------------------------------------------------------------------
# 3) All Trades sheet - close the open HighProbConvergence trade (row 139)
# -----------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(139, 7).Value = 0.95            # G - Exit Price
$allTrades.Cells.Item(139, 8).Value = "CLOSED"        # H - Status
$allTrades.Cells.Item(139, 9).Value = -1.0417         # I - P&L %
$allTrades.Cells.Item(139, 10).Value = -0.01          # J - P&L $
$allTrades.Cells.Item(139, 11).Value = 100.4          # K - Capital After
$allTrades.Cells.Item(139, 12).Value = "early_exit"   # L - Exit Reason
$allTrades.Cells.Item(139, 13).Value = 0.19           # M - Duration (min)

# New row 168: HighProbConvergence trade #167 (OPEN)
$allTrades.Cells.Item(168, 1).Value = 167
$allTrades.Cells.Item(168, 2).Value = "'2026-02-18"
$allTrades.Cells.Item(168, 2).Style = "Normal"
$allTrades.Cells.Item(168, 3).Value = "00:37:59"
$allTrades.Cells.Item(168, 4).Value = "HighProbConvergence"
$allTrades.Cells.Item(168, 5).Value = "UP"
$allTrades.Cells.Item(168, 6).Value = 0.96
$allTrades.Cells.Item(168, 8).Value = "OPEN"
$allTrades.Cells.Item(168, 9).Value = 0
$allTrades.Cells.Item(168, 10).Value = 0
$allTrades.Cells.Item(168, 11).Value = 100.4130057263667
$allTrades.Cells.Item(168, 13).Value = 0
$allTrades.Cells.Item(168, 14).Value = 0
$allTrades.Cells.Item(168, 15).Value = 0
$allTrades.Cells.Item(168, 16).Value = 0.95
$allTrades.Cells.Item(168, 17).Value = "Mean reversion UP: price 1.56% below mean (z=-2.00)"

# New row 169: MarketMaking trade #168 (OPEN)
$allTrades.Cells.Item(169, 1).Value = 168
$allTrades.Cells.Item(169, 2).Value = "'2026-02-18"
$allTrades.Cells.Item(169, 2).Style = "Normal"
$allTrades.Cells.Item(169, 3).Value = "00:37:59"
$allTrades.Cells.Item(169, 4).Value = "MarketMaking"
$allTrades.Cells.Item(169, 5).Value = "DOWN"
$allTrades.Cells.Item(169, 6).Value = 0.04
$allTrades.Cells.Item(169, 8).Value = "OPEN"
$allTrades.Cells.Item(169, 9).Value = 0
$allTrades.Cells.Item(169, 10).Value = 0
$allTrades.Cells.Item(169, 11).Value = 99.19858346467944
$allTrades.Cells.Item(169, 13).Value = 0
$allTrades.Cells.Item(169, 14).Value = 0
$allTrades.Cells.Item(169, 15).Value = 0
$allTrades.Cells.Item(169, 16).Value = 0.6
$allTrades.Cells.Item(169, 17).Value = "Normal spread capture: 198 bps"

# -----------------------------------------------------------------------
# 4) HighProbConvergence sheet - close trade (row 14) and append row 26
#    Column layout differs from "All Trades":
#    A..K same, L=Entry Slippage, M=Exit Slippage, N=Confidence,
#    O=Entry Reason, P=Exit Reason, Q=Duration (min)
# -----------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(14, 7).Value = 0.95             # G - Exit Price
$hpc.Cells.Item(14, 8).Value = "CLOSED"         # H - Status
$hpc.Cells.Item(14, 9).Value = -1.0417          # I - P&L %
$hpc.Cells.Item(14, 10).Value = -0.01           # J - P&L $
$hpc.Cells.Item(14, 11).Value = 100.4           # K - Capital After
$hpc.Cells.Item(14, 16).Value = "early_exit"    # P - Exit Reason
$hpc.Cells.Item(14, 17).Value = 0.19            # Q - Duration (min)

# New row 26: HighProbConvergence trade #167 (OPEN)
$hpc.Cells.Item(26, 1).Value = 167
$hpc.Cells.Item(26, 2).Value = "'2026-02-18"
$hpc.Cells.Item(26, 2).Style = "Normal"
$hpc.Cells.Item(26, 3).Value = "00:37:59"
$hpc.Cells.Item(26, 4).Value = "HighProbConvergence"
$hpc.Cells.Item(26, 5).Value = "UP"
$hpc.Cells.Item(26, 6).Value = 0.96
$hpc.Cells.Item(26, 8).Value = "OPEN"
$hpc.Cells.Item(26, 9).Value = 0
$hpc.Cells.Item(26, 10).Value = 0
$hpc.Cells.Item(26, 11).Value = 100.4130057263667
$hpc.Cells.Item(26, 12).Value = 0
$hpc.Cells.Item(26, 13).Value = 0
$hpc.Cells.Item(26, 14).Value = 0.95
$hpc.Cells.Item(26, 15).Value = "Mean reversion UP: price 1.56% below mean (z=-2.00)"
$hpc.Cells.Item(26, 17).Value = 0

# -----------------------------------------------------------------------
# 5) MarketMaking sheet - append row 66
# -----------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(66, 1).Value = 168
$mm.Cells.Item(66, 2).Value = "'2026-02-18"
$mm.Cells.Item(66, 2).Style = "Normal"
$mm.Cells.Item(66, 3).Value = "00:37:59"
$mm.Cells.Item(66, 4).Value = "MarketMaking"
$mm.Cells.Item(66, 5).Value = "DOWN"
$mm.Cells.Item(66, 6).Value = 0.04
$mm.Cells.Item(66, 8).Value = "OPEN"
$mm.Cells.Item(66, 9).Value = 0
$mm.Cells.Item(66, 10).Value = 0
$mm.Cells.Item(66, 11).Value = 99.19858346467944
$mm.Cells.Item(66, 12).Value = 0
$mm.Cells.Item(66, 13).Value = 0
$mm.Cells.Item(66, 14).Value = 0.6
$mm.Cells.Item(66, 15).Value = "Normal spread capture: 198 bps"
$mm.Cells.Item(66, 17).Value = 0
